$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-05 Friday" "2025-12-06 Saturday"

Replace-Text "580÷3=" "521÷3="
Replace-Text "384÷2=" "318÷9="
Replace-Text "960÷9=" "513÷3="
Replace-Text "365÷8=" "966÷5="
Replace-Text "767÷9=" "765÷2="
Replace-Text "926÷6=" "420÷9="
Replace-Text "770÷4=" "540÷3="
Replace-Text "467÷6=" "241÷6="
Replace-Text "962÷2=" "465÷4="
Replace-Text "395÷6=" "104÷4="
Replace-Text "574÷8=" "187÷7="
Replace-Text "594÷6=" "556÷2="
Replace-Text "297÷9=" "826÷3="
Replace-Text "897÷4=" "386÷3="
Replace-Text "164÷5=" "166÷6="
Replace-Text "466÷2=" "156÷6="
Replace-Text "671÷8=" "392÷2="
Replace-Text "225÷8=" "564÷5="
Replace-Text "360÷8=" "632÷2="
Replace-Text "387÷4=" "202÷7="
Replace-Text "356÷4=" "707÷6="
Replace-Text "573÷6=" "822÷3="
Replace-Text "278÷5=" "596÷8="
Replace-Text "192÷8=" "553÷9="
Replace-Text "866÷5=" "280÷5="
